$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and 1h-volume-change (E) columns with refreshed quotes.
# D-column values are entered via a temporary Text number format so that
# numeric-looking strings (e.g. "1.030") are kept as literal text, matching
# the source data; ClearFormats() afterwards restores the default (unstyled)
# cell format so no stray style is left behind.

$d = $ws.Range("D2")
$d.NumberFormat = '@'
$d.Value = '27.497.19'
$d.ClearFormats()
$ws.Range("E2").Value = '  +4.20%  '

$d = $ws.Range("D3")
$d.NumberFormat = '@'
$d.Value = '1.841.51'
$d.ClearFormats()
$ws.Range("E3").Value = '  +3.67%  '

$d = $ws.Range("D4")
$d.NumberFormat = '@'
$d.Value = '1.030'
$d.ClearFormats()
$ws.Range("E4").Value = '  +2.80%  '

$d = $ws.Range("D5")
$d.NumberFormat = '@'
$d.Value = '319.49'
$d.ClearFormats()
$ws.Range("E5").Value = '  +4.24%  '

$d = $ws.Range("D6")
$d.NumberFormat = '@'
$d.Value = '1.027'
$d.ClearFormats()
$ws.Range("E6").Value = '  +2.61%  '

$d = $ws.Range("D7")
$d.NumberFormat = '@'
$d.Value = '0.4369'
$d.ClearFormats()
$ws.Range("E7").Value = '  +3.16%  '

$d = $ws.Range("D8")
$d.NumberFormat = '@'
$d.Value = '0.3731'
$d.ClearFormats()
$ws.Range("E8").Value = '  +3.51%  '

$d = $ws.Range("D9")
$d.NumberFormat = '@'
$d.Value = '0.07388'
$d.ClearFormats()
$ws.Range("E9").Value = '  +3.29%  '

$d = $ws.Range("D10")
$d.NumberFormat = '@'
$d.Value = '0.8767'
$d.ClearFormats()
$ws.Range("E10").Value = '  +4.68%  '

$d = $ws.Range("D11")
$d.NumberFormat = '@'
$d.Value = '21.41'
$d.ClearFormats()
$ws.Range("E11").Value = '  +4.80%  '

$d = $ws.Range("D12")
$d.NumberFormat = '@'
$d.Value = '1.867.87'
$d.ClearFormats()
$ws.Range("E12").Value = '  +5.78%  '

$d = $ws.Range("D13")
$d.NumberFormat = '@'
$d.Value = '5.492'
$d.ClearFormats()
$ws.Range("E13").Value = '  +4.66%  '

$d = $ws.Range("D14")
$d.NumberFormat = '@'
$d.Value = '6.675'
$d.ClearFormats()
$ws.Range("E14").Value = '  +3.60%  '

$d = $ws.Range("D15")
$d.NumberFormat = '@'
$d.Value = '0.07147'
$d.ClearFormats()
$ws.Range("E15").Value = '  +3.69%  '

$d = $ws.Range("D16")
$d.NumberFormat = '@'
$d.Value = '82.65'
$d.ClearFormats()
$ws.Range("E16").Value = '  +4.52%  '

$ws.Range("E17").Value = '  +3.14%  '

$d = $ws.Range("D18")
$d.NumberFormat = '@'
$d.Value = '0.000009031'
$d.ClearFormats()
$ws.Range("E18").Value = '  +4.24%  '

$d = $ws.Range("D19")
$d.NumberFormat = '@'
$d.Value = '1.026'
$d.ClearFormats()
$ws.Range("E19").Value = '  +2.53%  '

$d = $ws.Range("D20")
$d.NumberFormat = '@'
$d.Value = '15.39'
$d.ClearFormats()
$ws.Range("E20").Value = '  +3.27%  '

$d = $ws.Range("D21")
$d.NumberFormat = '@'
$d.Value = '27.521.50'
$d.ClearFormats()
$ws.Range("E21").Value = '  +4.25%  '

$d = $ws.Range("D22")
$d.NumberFormat = '@'
$d.Value = '5.229'
$d.ClearFormats()
$ws.Range("E22").Value = '  +2.71%  '

$d = $ws.Range("D23")
$d.NumberFormat = '@'
$d.Value = '11.18'
$d.ClearFormats()
$ws.Range("E23").Value = '  +2.75%  '

$d = $ws.Range("D24")
$d.NumberFormat = '@'
$d.Value = '2.078.94'
$d.ClearFormats()
$ws.Range("E24").Value = '  +4.75%  '

$d = $ws.Range("D25")
$d.NumberFormat = '@'
$d.Value = '156.95'
$d.ClearFormats()
$ws.Range("E25").Value = '  +3.44%  '

$d = $ws.Range("D26")
$d.NumberFormat = '@'
$d.Value = '1.924'
$d.ClearFormats()
$ws.Range("E26").Value = '  +6.86%  '

$d = $ws.Range("D27")
$d.NumberFormat = '@'
$d.Value = '18.68'
$d.ClearFormats()
$ws.Range("E27").Value = '  +3.89%  '

$d = $ws.Range("D28")
$d.NumberFormat = '@'
$d.Value = '5.248'
$d.ClearFormats()
$ws.Range("E28").Value = '  +2.95%  '

$d = $ws.Range("D29")
$d.NumberFormat = '@'
$d.Value = '1.938'
$d.ClearFormats()
$ws.Range("E29").Value = '  +5.33%  '

$d = $ws.Range("D30")
$d.NumberFormat = '@'
$d.Value = '116.17'
$d.ClearFormats()
$ws.Range("E30").Value = '  +1.66%  '

$d = $ws.Range("D31")
$d.NumberFormat = '@'
$d.Value = '0.09065'
$d.ClearFormats()
$ws.Range("E31").Value = '  +2.69%  '

$d = $ws.Range("D32")
$d.NumberFormat = '@'
$d.Value = '1.209'
$d.ClearFormats()
$ws.Range("E32").Value = '  +7.53%  '

$d = $ws.Range("D33")
$d.NumberFormat = '@'
$d.Value = '0.7620'
$d.ClearFormats()
$ws.Range("E33").Value = '  +4.62%  '

$d = $ws.Range("D34")
$d.NumberFormat = '@'
$d.Value = '4.483'
$d.ClearFormats()
$ws.Range("E34").Value = '  +3.89%  '

$d = $ws.Range("D35")
$d.NumberFormat = '@'
$d.Value = '2.875'
$d.ClearFormats()
$ws.Range("E35").Value = '  +5.27%  '

$d = $ws.Range("D36")
$d.NumberFormat = '@'
$d.Value = '1.029'
$d.ClearFormats()
$ws.Range("E36").Value = '  +2.89%  '

$d = $ws.Range("D37")
$d.NumberFormat = '@'
$d.Value = '1.146'
$d.ClearFormats()
$ws.Range("E37").Value = '  +5.15%  '

$ws.Range("E38").Value = '  +4.48%  '

$d = $ws.Range("D39")
$d.NumberFormat = '@'
$d.Value = '0.05250'
$d.ClearFormats()
$ws.Range("E39").Value = '  +2.76%  '

$d = $ws.Range("D40")
$d.NumberFormat = '@'
$d.Value = '0.5171'
$d.ClearFormats()
$ws.Range("E40").Value = '  +5.29%  '

$d = $ws.Range("D41")
$d.NumberFormat = '@'
$d.Value = '2.782'
$d.ClearFormats()
$ws.Range("E41").Value = '  +6.99%  '

$d = $ws.Range("D42")
$d.NumberFormat = '@'
$d.Value = '0.1663'
$d.ClearFormats()
$ws.Range("E42").Value = '  +3.39%  '

$d = $ws.Range("D43")
$d.NumberFormat = '@'
$d.Value = '6.627'
$d.ClearFormats()
$ws.Range("E43").Value = '  +4.70%  '

$d = $ws.Range("D44")
$d.NumberFormat = '@'
$d.Value = '8.513'
$d.ClearFormats()
$ws.Range("E44").Value = '  +5.97%  '

$d = $ws.Range("D45")
$d.NumberFormat = '@'
$d.Value = '108.97'
$d.ClearFormats()
$ws.Range("E45").Value = '  +4.19%  '

$d = $ws.Range("D46")
$d.NumberFormat = '@'
$d.Value = '10.54'
$d.ClearFormats()
$ws.Range("E46").Value = '  +2.98%  '

$d = $ws.Range("D47")
$d.NumberFormat = '@'
$d.Value = '1.030'
$d.ClearFormats()
$ws.Range("E47").Value = '  +2.91%  '

$d = $ws.Range("D48")
$d.NumberFormat = '@'
$d.Value = '1.704'
$d.ClearFormats()
$ws.Range("E48").Value = '  +4.70%  '

$d = $ws.Range("D49")
$d.NumberFormat = '@'
$d.Value = '0.4635'
$d.ClearFormats()
$ws.Range("E49").Value = '  +4.27%  '

$d = $ws.Range("D50")
$d.NumberFormat = '@'
$d.Value = '0.06329'
$d.ClearFormats()
$ws.Range("E50").Value = '  +2.44%  '

$d = $ws.Range("D51")
$d.NumberFormat = '@'
$d.Value = '1.881'
$d.ClearFormats()
$ws.Range("E51").Value = '  +9.29%  '

